$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (column D) updates - values are stored as text, so force the cell
# number format to Text ("@") before assigning to avoid Excel auto-converting
# these numeric-looking strings into actual numbers.
$priceUpdates = @{
    "D2" = "248.31"
    "D3" = "21.69"
    "D4" = "5.288"
    "D6" = "3.429"
    "D7" = "6.382"
    "D8" = "0.8113"
    "D9" = "0.9498"
    "D10" = "0.1429"
    "D11" = "0.07609"
    "D12" = "0.03203"
    "D13" = "0.03094"
    "D14" = "0.09308"
    "D15" = "3.592"
    "D16" = "0.001597"
    "D17" = "0.04712"
    "D18" = "0.0005783"
    "D19" = "0.006249"
    "D20" = "0.005050"
    "D21" = "0.001036"
    "D23" = "3.786"
    "D25" = "0.3300"
    "D26" = "0.1302"
    "D40" = "0.03960"
    "D42" = "0.1063"
    "D43" = "0.003402"
    "D44" = "0.008801"
    "D45" = "0.00005612"
    "D47" = "0.0005503"
    "D48" = "0.7804"
    "D49" = "0.1754"
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
}

# Volume(1h) label (column E) updates - plain text, no special formatting needed.
$ws.Range("E18").Value = "17OneONE"
$ws.Range("E47").Value = "46ACDXExchangeACXTWorstin24h"
